$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: new "account reactivation" entry ---------------------------
$ws.Range("A8").Value = 44135.5704513889
$ws.Range("A8").NumberFormat = 'mm/dd/yy\ hh:mm\ AM/PM'

$ws.Range("B8").Value = "xxigua@example.com"
$ws.Range("C8").Value = "Xavier"
$ws.Range("D8").Value = "Xigua"
$ws.Range("E8").Value = "Mango"

$ws.Range("G8").Value = 44255
$ws.Range("G7").Copy()
$ws.Range("G8").PasteSpecial(-4122)

$ws.Range("H8").Value = "Power user"
$ws.Range("J8").Value = "Yes"

# Hyperlink for the new email address (mirrors B2..B7)
$ws.Hyperlinks.Add($ws.Range("B8"), "mailto:xxigua@example.com", "", "", "xxigua@example.com")

# Adding the hyperlink re-styles B8 with the built-in "Hyperlink" look;
# restore the plain formatting used by the rest of the column (matches
# B2/B6/B7, which are hyperlinked but carry the default style).
$ws.Range("B7").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$wb.Styles.Item("Hyperlink").Delete()

# --- Selection / view bookkeeping ---------------------------------------
$ws.Range("C11").Select()
